$d = $word.ActiveDocument
$tbl = $d.Tables(1)

function Set-RowText($rowIndex, $text) {
    $cell = $tbl.Rows($rowIndex).Cells(1)
    $rng = $cell.Range
    # Trim trailing cell-mark / paragraph-mark characters from the range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

# Simple single-value replacements (rows are 1-based in the Word OM)
Set-RowText 1 "0M"
Set-RowText 2 "0M"
Set-RowText 3 "0M"
Set-RowText 4 "33"
Set-RowText 5 "0.00002"
Set-RowText 6 "0.00077"
Set-RowText 7 "0.00021"
Set-RowText 8 "0.00006"
Set-RowText 9 "0.00045"
Set-RowText 10 "0.00047"
Set-RowText 11 "0.00050"
Set-RowText 12 "0.00892"

# Rows that collapse a multi-run/tab-delimited cell down to a single value
Set-RowText 44 "99.99"
Set-RowText 45 "0.01"
Set-RowText 46 "65"
